$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '56.801.34'
$ws.Cells.Item(2, 5).Value = '  -6.37%  '
$ws.Cells.Item(3, 4).Value = '2.612.46'
$ws.Cells.Item(3, 5).Value = '  -10.10%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '484.27'
$ws.Cells.Item(5, 5).Value = '  -8.77%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '134.38'
$ws.Cells.Item(6, 5).Value = '  -6.60%  '
$ws.Cells.Item(7, 5).Value = '  +0.22%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.505'
$ws.Cells.Item(8, 5).Value = '  -8.94%  '
$ws.Cells.Item(9, 4).Value = '2.604.71'
$ws.Cells.Item(9, 5).Value = '  -10.56%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '5.95'
$ws.Cells.Item(10, 5).Value = '  -0.89%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.100'
$ws.Cells.Item(11, 5).Value = '  -7.60%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.335'
$ws.Cells.Item(12, 5).Value = '  -7.82%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.125'
$ws.Cells.Item(13, 5).Value = '  -0.10%  '
$ws.Cells.Item(14, 4).Value = '3.145.24'
$ws.Cells.Item(14, 5).Value = '  -7.77%  '
$ws.Cells.Item(15, 4).Value = '57.018.99'
$ws.Cells.Item(15, 5).Value = '  -6.02%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '20.70'
$ws.Cells.Item(16, 5).Value = '  -8.78%  '
$ws.Cells.Item(17, 4).Value = '2.668.28'
$ws.Cells.Item(17, 5).Value = '  -8.37%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.0000130'
$ws.Cells.Item(18, 5).Value = '  -8.58%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.59'
$ws.Cells.Item(19, 5).Value = '  -8.87%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '10.56'
$ws.Cells.Item(20, 5).Value = '  -9.65%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '326.66'
$ws.Cells.Item(21, 5).Value = '  -10.87%  '
$ws.Cells.Item(22, 2).Value = 'Dai'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.999'
$ws.Cells.Item(22, 5).Value = '  -0.11%  '
$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.94'
$ws.Cells.Item(23, 5).Value = '  -10.22%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '5.62'
$ws.Cells.Item(24, 5).Value = '  -0.37%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '61.00'
$ws.Cells.Item(25, 5).Value = '  -5.72%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.412'
$ws.Cells.Item(26, 5).Value = '  -9.48%  '
$ws.Cells.Item(27, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.998'
$ws.Cells.Item(27, 5).Value = '  -0.12%  '
$ws.Cells.Item(28, 2).Value = 'Kaspa'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.162'
$ws.Cells.Item(28, 5).Value = '  -10.53%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.22'
$ws.Cells.Item(29, 5).Value = '  -7.77%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0791'
$ws.Cells.Item(30, 5).Value = '  -8.58%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.00'
$ws.Cells.Item(31, 5).Value = '  +0.04%  '
$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.55'
$ws.Cells.Item(32, 5).Value = '  -7.83%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '18.41'
$ws.Cells.Item(33, 5).Value = '  -6.69%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '146.42'
$ws.Cells.Item(34, 5).Value = '  -2.17%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.06'
$ws.Cells.Item(35, 5).Value = '  -7.29%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.25'
$ws.Cells.Item(36, 5).Value = '  -6.03%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.907'
$ws.Cells.Item(37, 5).Value = '  -9.55%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.08'
$ws.Cells.Item(38, 5).Value = '  -9.82%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '34.24'
$ws.Cells.Item(39, 5).Value = '  -8.95%  '
$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.997'
$ws.Cells.Item(40, 5).Value = '  +0.00%  '
$ws.Cells.Item(41, 2).Value = 'Filecoin'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.46'
$ws.Cells.Item(41, 5).Value = '  -6.81%  '
$ws.Cells.Item(42, 2).Value = 'Maker'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(42, 4).Value = '2.126.92'
$ws.Cells.Item(42, 5).Value = '  -7.14%  '
$ws.Cells.Item(43, 2).Value = 'Stacks'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.33'
$ws.Cells.Item(43, 5).Value = '  -10.77%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0544'
$ws.Cells.Item(44, 5).Value = '  -6.70%  '
$ws.Cells.Item(45, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '10.36'
$ws.Cells.Item(45, 5).Value = '  +0.34%  '
$ws.Cells.Item(46, 2).Value = 'Mantle'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.567'
$ws.Cells.Item(46, 5).Value = '  -12.36%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '18.19'
$ws.Cells.Item(47, 5).Value = '  -12.33%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0219'
$ws.Cells.Item(48, 5).Value = '  -7.16%  '
$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0862'
$ws.Cells.Item(49, 5).Value = '  -7.33%  '
$ws.Cells.Item(50, 2).Value = 'RenderToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '4.40'
$ws.Cells.Item(50, 5).Value = '  -10.97%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '17.38'
$ws.Cells.Item(51, 5).Value = '  -5.98%  '
